$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.995.91"
$ws.Range("E2").Value = "  -0.39%  "
$ws.Range("D3").Value = "2.053.73"
$ws.Range("E3").Value = "  +0.04%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'246.26"
$ws.Range("E5").Value = "  -1.40%  "
$ws.Range("E6").Value = "  -2.21%  "
$ws.Range("D7").Value = "'57.89"
$ws.Range("E7").Value = "  -4.06%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  -3.86%  "
$ws.Range("D10").Value = "'0.0783"
$ws.Range("E10").Value = "  -1.42%  "
$ws.Range("E11").Value = "  +2.11%  "
$ws.Range("D12").Value = "'15.30"
$ws.Range("E12").Value = "  -4.85%  "
$ws.Range("E13").Value = "  +5.46%  "
$ws.Range("E14").Value = "  +0.05%  "
$ws.Range("E15").Value = "  -2.64%  "
$ws.Range("D16").Value = "2.050.12"
$ws.Range("E16").Value = "  -0.06%  "
$ws.Range("D17").Value = "'17.91"
$ws.Range("E17").Value = "  -2.76%  "
$ws.Range("D18").Value = "36.946.28"
$ws.Range("E18").Value = "  -0.44%  "
$ws.Range("D19").Value = "'73.66"
$ws.Range("E19").Value = "  -3.66%  "
$ws.Range("D20").Value = "0.0₃0896"
$ws.Range("E20").Value = "  -1.29%  "
$ws.Range("D21").Value = "'5.43"
$ws.Range("E21").Value = "  -0.15%  "
$ws.Range("D22").Value = "'235.79"
$ws.Range("E22").Value = "  -1.24%  "
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("E24").Value = "  +1.86%  "
$ws.Range("D25").Value = "'10.53"
$ws.Range("E25").Value = "  +11.40%  "
$ws.Range("D26").Value = "'2.18"
$ws.Range("E26").Value = "  -2.13%  "
$ws.Range("D27").Value = "'169.34"
$ws.Range("E27").Value = "  -0.07%  "
$ws.Range("E28").Value = "  -1.24%  "
$ws.Range("D29").Value = "'5.56"
$ws.Range("E29").Value = "  +14.99%  "
$ws.Range("E30").Value = "  -2.04%  "
$ws.Range("E31").Value = "  -2.58%  "
$ws.Range("D32").Value = "'4.72"
$ws.Range("E32").Value = "  +3.59%  "
$ws.Range("D33").Value = "'0.0617"
$ws.Range("E33").Value = "  -2.15%  "
$ws.Range("D34").Value = "'2.35"
$ws.Range("E34").Value = "  +4.49%  "
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("D36").Value = "'1.83"
$ws.Range("E36").Value = "  +5.09%  "
$ws.Range("D37").Value = "'0.0818"
$ws.Range("E37").Value = "  -7.63%  "
$ws.Range("E38").Value = "  -1.52%  "
$ws.Range("D39").Value = "'5.17"
$ws.Range("E39").Value = "  -0.16%  "
$ws.Range("E40").Value = "  -2.99%  "
$ws.Range("E41").Value = "  -0.65%  "
$ws.Range("E42").Value = "  +0.83%  "
$ws.Range("D43").Value = "'0.0958"
$ws.Range("E43").Value = "  -11.05%  "
$ws.Range("D44").Value = "'97.15"
$ws.Range("E44").Value = "  -0.26%  "
$ws.Range("D45").Value = "'16.93"
$ws.Range("E45").Value = "  -4.84%  "
$ws.Range("D46").Value = "1.308.29"
$ws.Range("E46").Value = "  +1.17%  "
$ws.Range("E47").Value = "  -6.14%  "
$ws.Range("E48").Value = "  -0.31%  "
$ws.Range("E49").Value = "  -1.81%  "
$ws.Range("D50").Value = "2.236.49"
$ws.Range("E50").Value = "  -0.19%  "
$ws.Range("D51").Value = "'44.76"
$ws.Range("E51").Value = "  -0.45%  "
